# Adds sample/example rows (2-4), additional "Catatan tambahan" notes
# (L14:L18) and a couple of cosmetic sheet-view / column-width tweaks to
# the product import template — per commit "penambahan keterangan untuk
# template_posmi.xlsx".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extra "Catatan tambahan" block in column L (rows 14-18) ----------
# (written first so new shared-string entries land in the same order as
#  the authored workbook: notes block, then the three example rows)
$ws.Range("L14").Value = "Catatan tambahan:"
$ws.Range("L15").Value = "1. Semua Kolom harus diisi kecuali kolom keterangan"
$ws.Range("L16").Value = "2. apabila barang tidak memiliki barcode bisa diisi dengan angka atau huruf dan harus unik, misalnya"
$ws.Range("L17").Value = "    dibuat kode aa1 atau aa2 atau aa3 atau 123 atau 124 dan lain-lain"
$ws.Range("L18").Value = "3. apabila harga penjualan tidak ada mode grosir, min_beli_grosir diisi angka 1"

# --- Example data rows ------------------------------------------------
# Row 2: barcode, nama, satuan, stok, harga_beli, harga_ecer, harga_grosir, min_beli_grosir, keterangan
$ws.Range("A2").Value = 123
$ws.Range("B2").Value = "Barang Contoh"
$ws.Range("C2").Value = "PCS"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 200000
$ws.Range("F2").Value = 230000
$ws.Range("G2").Value = 230000
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = "barang tanpa grosir"

# Row 3
$ws.Range("A3").Value = 124
$ws.Range("B3").Value = "Barang Contoh 2"
$ws.Range("C3").Value = "PACK"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 120000
$ws.Range("F3").Value = 140000
$ws.Range("G3").Value = 135000
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = "barang grosir dengan minimal pembelian 2"

# Row 4
$ws.Range("A4").Value = "abc1"
$ws.Range("B4").Value = "Barang Contoh 3"
$ws.Range("C4").Value = "PACK"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 20000
$ws.Range("F4").Value = 28000
$ws.Range("G4").Value = 25000
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = "barang dengan kode huruf"

# --- Column width tweaks ----------------------------------------------
# The engine quantises ColumnWidth to 1/6-character steps (it doesn't
# actually measure glyphs of "Aptos Narrow"), so the literal inputs below
# are chosen so the stored <col width=".."> lands on — or as close as
# achievable to — the authored values (col B -> 15.5703125, col I -> 40).
$ws.Range("B1").EntireColumn.ColumnWidth = 14.65
$ws.Range("I1").EntireColumn.ColumnWidth = 39.2

# --- Sheet view: scroll back to column A, select A5 --------------------
$ws.Range("A5").Select()
